# Applies the "updated data files to include princeton public library" commit.
# - Renames the preserveMtLakes1000 entry (row 10) to preserveMtLakes, including
#   its image URL, and refreshes its recomputed RGB/cluster statistics.
# - Refreshes recomputed RGB/cluster statistics on several other rows (2,4,5,8,14,16,17)
#   that shifted slightly once the new image was added to the underlying analysis.
# - Appends a new row (20) for the Princeton Public Library (urbanLibrary) satellite image.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Row 2
$ws.Cells.Item(2, 17).Value2 = 0.204951708638044  # Q2
$ws.Cells.Item(2, 18).Value2 = 0.361271490115523  # R2
$ws.Cells.Item(2, 19).Value2 = 0.103829285737177  # S2
$ws.Cells.Item(2, 20).Value2 = 0.602941124851978  # T2
$ws.Cells.Item(2, 21).Value2 = 0.312061512218915  # U2
$ws.Cells.Item(2, 22).Value2 = 0.454567598926795  # V2
$ws.Cells.Item(2, 23).Value2 = 0.176471605660498  # W2
$ws.Cells.Item(2, 24).Value2 = 0.584050677081246  # X2
$ws.Cells.Item(2, 25).Value2 = 0.191569882362021  # Y2
$ws.Cells.Item(2, 26).Value2 = 0.272835318092153  # Z2
$ws.Cells.Item(2, 27).Value2 = 0.128496372697394  # AA2
$ws.Cells.Item(2, 28).Value2 = 0.515793847681763  # AB2
$ws.Cells.Item(2, 29).Value2 = 0.317409688390477  # AC2
$ws.Cells.Item(2, 30).Value2 = 0.315146612740142  # AD2
$ws.Cells.Item(2, 31).Value2 = 0.297610074455373  # AE2
$ws.Cells.Item(2, 32).Value2 = 0.0698336244140086  # AF2

# Row 4
$ws.Cells.Item(4, 17).Value2 = 0.341761781187704  # Q4
$ws.Cells.Item(4, 18).Value2 = 0.193506370174901  # R4
$ws.Cells.Item(4, 19).Value2 = 0.0985673297631171  # S4
$ws.Cells.Item(4, 20).Value2 = 0.551805091994163  # T4
$ws.Cells.Item(4, 21).Value2 = 0.441594237333882  # U4
$ws.Cells.Item(4, 22).Value2 = 0.301944337935679  # V4
$ws.Cells.Item(4, 23).Value2 = 0.171100168762641  # W4
$ws.Cells.Item(4, 24).Value2 = 0.542373134232376  # X4
$ws.Cells.Item(4, 25).Value2 = 0.258447097760468  # Y4
$ws.Cells.Item(4, 26).Value2 = 0.184214404599468  # Z4
$ws.Cells.Item(4, 27).Value2 = 0.122766917951963  # AA4
$ws.Cells.Item(4, 28).Value2 = 0.476703387160047  # AB4
$ws.Cells.Item(4, 29).Value2 = 0.407344425039066  # AC4
$ws.Cells.Item(4, 30).Value2 = 0.30290651714312  # AD4
$ws.Cells.Item(4, 31).Value2 = 0.233696111774979  # AE4
$ws.Cells.Item(4, 32).Value2 = 0.0560529460428348  # AF4

# Row 5
$ws.Cells.Item(5, 17).Value2 = 0.199146374695632  # Q5
$ws.Cells.Item(5, 18).Value2 = 0.0986747809492752  # R5
$ws.Cells.Item(5, 19).Value2 = 0.370998162528187  # S5
$ws.Cells.Item(5, 20).Value2 = 0.589293581623865  # T5
$ws.Cells.Item(5, 21).Value2 = 0.293962030251813  # U5
$ws.Cells.Item(5, 22).Value2 = 0.154170968148882  # V5
$ws.Cells.Item(5, 23).Value2 = 0.444230737084637  # W5
$ws.Cells.Item(5, 24).Value2 = 0.564326079821592  # X5
$ws.Cells.Item(5, 25).Value2 = 0.180585158981223  # Y5
$ws.Cells.Item(5, 26).Value2 = 0.108536048177919  # Z5
$ws.Cells.Item(5, 27).Value2 = 0.278994443865823  # AA5
$ws.Cells.Item(5, 28).Value2 = 0.528382141541654  # AB5
$ws.Cells.Item(5, 29).Value2 = 0.402266752458866  # AC5
$ws.Cells.Item(5, 30).Value2 = 0.323160216931703  # AD5
$ws.Cells.Item(5, 31).Value2 = 0.224270613107822  # AE5
$ws.Cells.Item(5, 32).Value2 = 0.0503024175016086  # AF5

# Row 8
$ws.Cells.Item(8, 17).Value2 = 0.223813620423603  # Q8
$ws.Cells.Item(8, 18).Value2 = 0.115189278635276  # R8
$ws.Cells.Item(8, 19).Value2 = 0.40345144898967  # S8
$ws.Cells.Item(8, 20).Value2 = 0.586733777868203  # T8
$ws.Cells.Item(8, 21).Value2 = 0.328427491171107  # U8
$ws.Cells.Item(8, 22).Value2 = 0.194392081942834  # V8
$ws.Cells.Item(8, 23).Value2 = 0.465319231375064  # W8
$ws.Cells.Item(8, 24).Value2 = 0.568000171052044  # X8
$ws.Cells.Item(8, 25).Value2 = 0.218584157153304  # Y8
$ws.Cells.Item(8, 26).Value2 = 0.143660635614015  # Z8
$ws.Cells.Item(8, 27).Value2 = 0.320048403622281  # AA8
$ws.Cells.Item(8, 28).Value2 = 0.496879536261878  # AB8
$ws.Cells.Item(8, 29).Value2 = 0.361272175751448  # AC8
$ws.Cells.Item(8, 30).Value2 = 0.265590587370163  # AD8
$ws.Cells.Item(8, 31).Value2 = 0.265239452155529  # AE8
$ws.Cells.Item(8, 32).Value2 = 0.107897784722861  # AF8

# Row 10
$ws.Cells.Item(10, 1).Value2 = 'preserveMtLakes'  # A10
$ws.Cells.Item(10, 2).Value2 = 'https://github.com/Imageomics/Andromeda/raw/main/datasets/satelliteData/preserveMtLakes.png'  # B10
$ws.Cells.Item(10, 11).Value2 = 0.274509803921569  # K10
$ws.Cells.Item(10, 12).Value2 = 0.162573166713187  # L10
$ws.Cells.Item(10, 13).Value2 = 0.364705882352941  # M10
$ws.Cells.Item(10, 14).Value2 = 0.141980590942576  # N10
$ws.Cells.Item(10, 15).Value2 = 0.294117647058824  # O10
$ws.Cells.Item(10, 16).Value2 = 0.14210989979474  # P10
$ws.Cells.Item(10, 17).Value2 = 0.267661562189387  # Q10
$ws.Cells.Item(10, 18).Value2 = 0.150240993285796  # R10
$ws.Cells.Item(10, 19).Value2 = 0.437783344077499  # S10
$ws.Cells.Item(10, 20).Value2 = 0.632509539636201  # T10
$ws.Cells.Item(10, 21).Value2 = 0.356086932939028  # U10
$ws.Cells.Item(10, 22).Value2 = 0.21820985583095  # V10
$ws.Cells.Item(10, 23).Value2 = 0.481264137518036  # W10
$ws.Cells.Item(10, 24).Value2 = 0.627157761440502  # X10
$ws.Cells.Item(10, 25).Value2 = 0.289000651275803  # Y10
$ws.Cells.Item(10, 26).Value2 = 0.208345511780479  # Z10
$ws.Cells.Item(10, 27).Value2 = 0.416060480132297  # AA10
$ws.Cells.Item(10, 28).Value2 = 0.635005783590392  # AB10
$ws.Cells.Item(10, 29).Value2 = 0.366849894291755  # AC10
$ws.Cells.Item(10, 30).Value2 = 0.282005699053222  # AD10
$ws.Cells.Item(10, 31).Value2 = 0.236027208383123  # AE10
$ws.Cells.Item(10, 32).Value2 = 0.1151171982719  # AF10
$ws.Cells.Item(10, 33).Value2 = 'gray30'  # AG10
$ws.Cells.Item(10, 34).Value2 = 'gray19'  # AH10
$ws.Cells.Item(10, 35).Value2 = 'gray45'  # AI10

# Row 14
$ws.Cells.Item(14, 17).Value2 = 0.415040407962094  # Q14
$ws.Cells.Item(14, 18).Value2 = 0.233873944664757  # R14
$ws.Cells.Item(14, 19).Value2 = 0.108223305346893  # S14
$ws.Cells.Item(14, 20).Value2 = 0.579102404647386  # T14
$ws.Cells.Item(14, 21).Value2 = 0.441550101742689  # U14
$ws.Cells.Item(14, 22).Value2 = 0.323920710338058  # V14
$ws.Cells.Item(14, 23).Value2 = 0.187500121892146  # W14
$ws.Cells.Item(14, 24).Value2 = 0.561522360690173  # X14
$ws.Cells.Item(14, 25).Value2 = 0.329398239396863  # Y14
$ws.Cells.Item(14, 26).Value2 = 0.218969231435148  # Z14
$ws.Cells.Item(14, 27).Value2 = 0.144009569346886  # AA14
$ws.Cells.Item(14, 28).Value2 = 0.524350564735912  # AB14
$ws.Cells.Item(14, 29).Value2 = 0.327281919294053  # AC14
$ws.Cells.Item(14, 30).Value2 = 0.292311793363361  # AD14
$ws.Cells.Item(14, 31).Value2 = 0.266156815883813  # AE14
$ws.Cells.Item(14, 32).Value2 = 0.114249471458774  # AF14

# Row 16
$ws.Cells.Item(16, 17).Value2 = 0.288346720349506  # Q16
$ws.Cells.Item(16, 18).Value2 = 0.168390126059778  # R16
$ws.Cells.Item(16, 19).Value2 = 0.454040793920364  # S16
$ws.Cells.Item(16, 20).Value2 = 0.649498305401255  # T16
$ws.Cells.Item(16, 21).Value2 = 0.374038405313683  # U16
$ws.Cells.Item(16, 22).Value2 = 0.227670648318794  # V16
$ws.Cells.Item(16, 23).Value2 = 0.490080066013682  # W16
$ws.Cells.Item(16, 24).Value2 = 0.637245095414923  # X16
$ws.Cells.Item(16, 25).Value2 = 0.305989040406259  # Y16
$ws.Cells.Item(16, 26).Value2 = 0.233742913421416  # Z16
$ws.Cells.Item(16, 27).Value2 = 0.435221311859108  # AA16
$ws.Cells.Item(16, 28).Value2 = 0.629861898972686  # AB16
$ws.Cells.Item(16, 29).Value2 = 0.331659159849251  # AC16
$ws.Cells.Item(16, 30).Value2 = 0.289896130159022  # AD16
$ws.Cells.Item(16, 31).Value2 = 0.247952936850813  # AE16
$ws.Cells.Item(16, 32).Value2 = 0.130491773140914  # AF16

# Row 17
$ws.Cells.Item(17, 17).Value2 = 0.292485254009982  # Q17
$ws.Cells.Item(17, 18).Value2 = 0.170953388221756  # R17
$ws.Cells.Item(17, 19).Value2 = 0.45861257162221  # S17
$ws.Cells.Item(17, 21).Value2 = 0.381564529542702  # U17
$ws.Cells.Item(17, 22).Value2 = 0.227219813435469  # V17
$ws.Cells.Item(17, 23).Value2 = 0.494421844158163  # W17
$ws.Cells.Item(17, 25).Value2 = 0.309013968191812  # Y17
$ws.Cells.Item(17, 26).Value2 = 0.237773962004343  # Z17
$ws.Cells.Item(17, 27).Value2 = 0.436045381692298  # AA17
$ws.Cells.Item(17, 29).Value2 = 0.322654655758801  # AC17
$ws.Cells.Item(17, 30).Value2 = 0.282461623310966  # AD17
$ws.Cells.Item(17, 31).Value2 = 0.257404173177682  # AE17

# Row 20
$ws.Cells.Item(20, 1).Value2 = 'urbanLibrary'  # A20
$ws.Cells.Item(20, 2).Value2 = 'https://github.com/Imageomics/Andromeda/raw/main/datasets/satelliteData/urbanLibrary.png'  # B20
$ws.Cells.Item(20, 3).Value2 = 'Princeton Public Library'  # C20
$ws.Cells.Item(20, 4).Value2 = 'urban'  # D20
$ws.Cells.Item(20, 5).Value2 = 40.35151  # E20
$ws.Cells.Item(20, 6).Value2 = -74.660325  # F20
$ws.Cells.Item(20, 7).Value2 = 40.357259  # G20
$ws.Cells.Item(20, 8).Value2 = -74.671573  # H20
$ws.Cells.Item(20, 9).Value2 = 40.34662  # I20
$ws.Cells.Item(20, 10).Value2 = -74.648996  # J20
$ws.Cells.Item(20, 11).Value2 = 0.298039215686275  # K20
$ws.Cells.Item(20, 12).Value2 = 0.163940995974375  # L20
$ws.Cells.Item(20, 13).Value2 = 0.380392156862745  # M20
$ws.Cells.Item(20, 14).Value2 = 0.144352586153883  # N20
$ws.Cells.Item(20, 15).Value2 = 0.313725490196078  # O20
$ws.Cells.Item(20, 16).Value2 = 0.141369147301135  # P20
$ws.Cells.Item(20, 17).Value2 = 0.290054309778663  # Q20
$ws.Cells.Item(20, 18).Value2 = 0.167018335101561  # R20
$ws.Cells.Item(20, 19).Value2 = 0.458670411202548  # S20
$ws.Cells.Item(20, 20).Value2 = 0.643780422542729  # T20
$ws.Cells.Item(20, 21).Value2 = 0.373545394552962  # U20
$ws.Cells.Item(20, 22).Value2 = 0.229170677153524  # V20
$ws.Cells.Item(20, 23).Value2 = 0.493155287112214  # W20
$ws.Cells.Item(20, 24).Value2 = 0.63821369359121  # X20
$ws.Cells.Item(20, 25).Value2 = 0.309083010639517  # Y20
$ws.Cells.Item(20, 26).Value2 = 0.230508165471853  # Z20
$ws.Cells.Item(20, 27).Value2 = 0.439155055375326  # AA20
$ws.Cells.Item(20, 28).Value2 = 0.648456573677854  # AB20
$ws.Cells.Item(20, 29).Value2 = 0.34112510341024  # AC20
$ws.Cells.Item(20, 30).Value2 = 0.298858350951374  # AD20
$ws.Cells.Item(20, 31).Value2 = 0.243283390017465  # AE20
$ws.Cells.Item(20, 32).Value2 = 0.116733155620921  # AF20
$ws.Cells.Item(20, 33).Value2 = 'gray32'  # AG20
$ws.Cells.Item(20, 34).Value2 = 'gray21'  # AH20
$ws.Cells.Item(20, 35).Value2 = 'gray46'  # AI20
$ws.Cells.Item(20, 36).Value2 = 'gray64'  # AJ20
